# Applies the "INFOTILES" fixes:
#  - TABLE_FIELDS_config!B3: "infotile" -> "infotiles" (table database name typo fix)
#  - TABLE_FIELDS_config!A4: "localname" -> "Display name"
#  - TABLE_FIELDS_config!E11: clarify the department/office relation comment
#  - Make TABLE_FIELDS_config the active/selected sheet (was LISA_screenshots)
#  - Update the remembered selection on TABLE_FIELDS_config to B13

$wb = $excel.ActiveWorkbook

$wsConfig = $wb.Worksheets.Item("TABLE_FIELDS_config")

$wsConfig.Range("B3").Value = "infotiles"
$wsConfig.Range("A4").Value = "Display name"
$wsConfig.Range("E11").Value = "Relation to [department] or office. NOT MANDATORY - functionality can be turned of in VBA"

# Move the active tab back to TABLE_FIELDS_config, with the remembered
# selection on B13, and make sure LISA_screenshots is no longer the
# "tabSelected" sheet.
$wsConfig.Activate()
$wsConfig.Range("B13").Select()
